# Insert a new data row at row 313 (pushing existing rows 313:408 down to 314:409),
# matching the author's commit: a new weekly "Fruta / hortaliza" price observation
# was inserted into the historical series for Piña (Macroferia Regional de Talca).
#
# The new row duplicates the values of what was previously the last row in the
# sheet (old row 408), except for the date, which is the new, more recent
# observation date (serial 45093 = 2023-06-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 313; existing rows 313:408 shift to 314:409.
$ws.Rows.Item(313).Insert()

# Populate the newly inserted row 313 with the new observation.
$ws.Range("A313").Value2 = 5
$ws.Range("B313").Value = "Macroferia Regional de Talca"
$ws.Range("C313").Value = "Maule"
$ws.Range("D313").Value2 = 45093
$ws.Range("E313").Value2 = 7
$ws.Range("F313").Value = "Fruta"
$ws.Range("G313").Value2 = 100108
$ws.Range("H313").Value = "Tropicales y subtropicales"
$ws.Range("I313").Value2 = 100108005
$ws.Range("J313").Value = "Piña"
$ws.Range("K313").Value = "Caramelo"
$ws.Range("L313").Value = "Segunda"
$ws.Range("M313").Value2 = 250
$ws.Range("N313").Value2 = 21000
$ws.Range("O313").Value2 = 21000
$ws.Range("P313").Value2 = 21000
$ws.Range("Q313").Value = "`$/caja 14 unidades"
$ws.Range("R313").Value = "Ecuador"
$ws.Range("S313").Value2 = 1500
$ws.Range("T313").Value2 = 14
